$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 5
$ws.Range("G3").Value = -3
$ws.Range("H3").Value = 13

$ws.Range("D3").Select()
